$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "vote16" row (row 19: vote16 / "remember if voted in 2016 election")
# This shifts all subsequent rows up by one, removing the now-empty last row (29).
$ws.Rows("19:19").Delete()

# Update the active cell selection as recorded in the saved file.
$ws.Range("E9").Select()
